$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1996434937611408
$ws.Range("C2").Value = 0.5543672014260249
$ws.Range("J2").Value = 0.0053475935828877
$ws.Range("P2").Value = 0.1515151515151515
$ws.Range("S2").Value = 0.08912655971479501
$ws.Range("B3").Value = 0.00308641975308642
$ws.Range("C3").Value = 0.02777777777777778
$ws.Range("J3").Value = 0.02469135802469136
$ws.Range("P3").Value = 0.7098765432098766
$ws.Range("S3").Value = 0.2345679012345679
$ws.Range("J4").Value = 0.075
$ws.Range("O4").Value = 0.0125
$ws.Range("P4").Value = 0.6375
$ws.Range("S4").Value = 0.275
$ws.Range("B6").Value = 0.06478873239436619
$ws.Range("D6").Value = 0.01971830985915493
$ws.Range("F6").Value = 0.04788732394366197
$ws.Range("J6").Value = 0.2788732394366197
$ws.Range("O6").Value = 0.008450704225352112
$ws.Range("Q6").Value = 0.1380281690140845
$ws.Range("R6").Value = 0.05352112676056338
$ws.Range("S6").Value = 0.3887323943661972
$ws.Range("B7").Value = 0.09641873278236915
$ws.Range("D7").Value = 0.01928374655647383
$ws.Range("E7").Value = 0.002754820936639119
$ws.Range("F7").Value = 0.03305785123966942
$ws.Range("J7").Value = 0.1680440771349862
$ws.Range("O7").Value = 0.01652892561983471
$ws.Range("Q7").Value = 0.162534435261708
$ws.Range("R7").Value = 0.07162534435261708
$ws.Range("S7").Value = 0.4297520661157025
$ws.Range("B8").Value = 0.09655172413793103
$ws.Range("D8").Value = 0.01839080459770115
$ws.Range("F8").Value = 0.05172413793103448
$ws.Range("J8").Value = 0.1252873563218391
$ws.Range("O8").Value = 0.009195402298850575
$ws.Range("Q8").Value = 0.1896551724137931
$ws.Range("R8").Value = 0.06666666666666667
$ws.Range("S8").Value = 0.4425287356321839
$ws.Range("B9").Value = 0.09510086455331412
$ws.Range("D9").Value = 0.02305475504322766
$ws.Range("F9").Value = 0.08357348703170028
$ws.Range("J9").Value = 0.1037463976945245
$ws.Range("O9").Value = 0.01152737752161383
$ws.Range("Q9").Value = 0.1873198847262248
$ws.Range("R9").Value = 0.08645533141210375
$ws.Range("S9").Value = 0.4092219020172911
$ws.Range("B10").Value = 0.1184560780834073
$ws.Range("D10").Value = 0.02129547471162378
$ws.Range("E10").Value = 0.0004436557231588288
$ws.Range("F10").Value = 0.05989352262644188
$ws.Range("J10").Value = 0.1206743566992014
$ws.Range("O10").Value = 0.01552795031055901
$ws.Range("Q10").Value = 0.2076308784383319
$ws.Range("R10").Value = 0.064773735581189
$ws.Range("S10").Value = 0.391304347826087
$ws.Range("G11").Value = 0.137221269296741
$ws.Range("J11").Value = 0.09433962264150944
$ws.Range("K11").Value = 0.1989708404802744
$ws.Range("L11").Value = 0.5540308747855918
$ws.Range("S11").Value = 0.01543739279588336
$ws.Range("G12").Value = 0.7329376854599406
$ws.Range("J12").Value = 0.2047477744807122
$ws.Range("K12").Value = 0.008902077151335312
$ws.Range("L12").Value = 0.02373887240356083
$ws.Range("S12").Value = 0.02967359050445104
$ws.Range("G13").Value = 0.6865671641791045
$ws.Range("J13").Value = 0.2985074626865671
$ws.Range("S13").Value = 0.01492537313432836
$ws.Range("F15").Value = 0.01842105263157895
$ws.Range("H15").Value = 0.1605263157894737
$ws.Range("I15").Value = 0.07368421052631578
$ws.Range("J15").Value = 0.35
$ws.Range("K15").Value = 0.07105263157894737
$ws.Range("M15").Value = 0.007894736842105263
$ws.Range("O15").Value = 0.05
$ws.Range("S15").Value = 0.2684210526315789
$ws.Range("F16").Value = 0.01396648044692737
$ws.Range("H16").Value = 0.1731843575418995
$ws.Range("I16").Value = 0.07541899441340782
$ws.Range("J16").Value = 0.3854748603351955
$ws.Range("K16").Value = 0.09497206703910614
$ws.Range("M16").Value = 0.0223463687150838
$ws.Range("O16").Value = 0.04748603351955307
$ws.Range("S16").Value = 0.1871508379888268
$ws.Range("F17").Value = 0.01496259351620948
$ws.Range("H17").Value = 0.1770573566084788
$ws.Range("I17").Value = 0.09351620947630923
$ws.Range("J17").Value = 0.3753117206982544
$ws.Range("K17").Value = 0.1184538653366584
$ws.Range("M17").Value = 0.01870324189526185
$ws.Range("O17").Value = 0.0885286783042394
$ws.Range("S17").Value = 0.1134663341645885
$ws.Range("F18").Value = 0.01444043321299639
$ws.Range("H18").Value = 0.1985559566787004
$ws.Range("I18").Value = 0.05776173285198556
$ws.Range("J18").Value = 0.444043321299639
$ws.Range("K18").Value = 0.08664259927797834
$ws.Range("M18").Value = 0.007220216606498195
$ws.Range("N18").Value = 0.003610108303249098
$ws.Range("O18").Value = 0.05776173285198556
$ws.Range("S18").Value = 0.1299638989169675
$ws.Range("F19").Value = 0.0162999185004075
$ws.Range("H19").Value = 0.2249388753056235
$ws.Range("I19").Value = 0.08353708231458842
$ws.Range("J19").Value = 0.3492257538712307
$ws.Range("K19").Value = 0.1145069274653627
$ws.Range("M19").Value = 0.01711491442542787
$ws.Range("N19").Value = 0.00162999185004075
$ws.Range("O19").Value = 0.06316218418907905
$ws.Range("S19").Value = 0.1295843520782396
